$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New file handed off: e230a072-2246-4ea1-abee-1f9de0b94b68.md
# It takes the row that used to hold ".localization-config" (row 7) on every
# sheet, and ".localization-config" is pushed down to a new row 8.
# ---------------------------------------------------------------------------

$mdName   = "e230a072-2246-4ea1-abee-1f9de0b94b68.md"
$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/1b1c9ff59f406b66eb32bf633a2187f15b6eb164/e2e/e230a072-2246-4ea1-abee-1f9de0b94b68.md"
$cfgName  = ".localization-config"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/1b1c9ff59f406b66eb32bf633a2187f15b6eb164/.localization-config"

# ----------------------------- Overview sheet ------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(8, 1).Value = $cfgName
$ov.Cells.Item(8, 2).Value = "Not to be localized"
$ov.Cells.Item(8, 3).Value = "Not to be localized"

$ov.Cells.Item(7, 1).Value = $mdName
$ov.Cells.Item(7, 2).Value = "Ready for handoff"
$ov.Cells.Item(7, 3).Value = "Ready for handoff"

$ov.Hyperlinks.Add($ov.Cells.Item(7, 1), $mdUrl, "", "", $mdName) | Out-Null
$ov.Hyperlinks.Add($ov.Cells.Item(8, 1), $cfgUrl, "", "", $cfgName) | Out-Null

# ------------------------------ zh-cn sheet --------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhXlfName = "e230a072-2246-4ea1-abee-1f9de0b94b68.a3c2af93fd0cc01fae1536b1552fde9454ae458b.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3c2af93fd0cc01fae1536b1552fde9454ae458b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlfName"

$zh.Cells.Item(8, 1).Value = $cfgName
$zh.Cells.Item(8, 2).Value = "Not to be localized"
$zh.Cells.Item(8, 4).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(8, 7).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(8, 8).Value = "Ignored"

$zh.Cells.Item(7, 1).Value = $mdName
$zh.Cells.Item(7, 2).Value = "Ready for handoff"
$zh.Cells.Item(7, 3).Value = $zhXlfName
$zh.Cells.Item(7, 4).Value = "2016-03-07 04:22:59"
$zh.Cells.Item(7, 7).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(7, 8).Value = "Include"

$zh.Hyperlinks.Add($zh.Cells.Item(7, 1), $mdUrl, "", "", $mdName) | Out-Null
$zh.Hyperlinks.Add($zh.Cells.Item(7, 3), $zhXlfUrl, "", "", $zhXlfName) | Out-Null
$zh.Hyperlinks.Add($zh.Cells.Item(8, 1), $cfgUrl, "", "", $cfgName) | Out-Null

# ------------------------------ de-de sheet --------------------------------
$de = $wb.Worksheets.Item("de-de")

$deXlfName = "e230a072-2246-4ea1-abee-1f9de0b94b68.a3c2af93fd0cc01fae1536b1552fde9454ae458b.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3c2af93fd0cc01fae1536b1552fde9454ae458b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlfName"

$de.Cells.Item(8, 1).Value = $cfgName
$de.Cells.Item(8, 2).Value = "Not to be localized"
$de.Cells.Item(8, 4).Value = "0001-01-01 00:00:00"
$de.Cells.Item(8, 7).Value = "0001-01-01 00:00:00"
$de.Cells.Item(8, 8).Value = "Ignored"

$de.Cells.Item(7, 1).Value = $mdName
$de.Cells.Item(7, 2).Value = "Ready for handoff"
$de.Cells.Item(7, 3).Value = $deXlfName
$de.Cells.Item(7, 4).Value = "2016-03-07 04:23:09"
$de.Cells.Item(7, 7).Value = "0001-01-01 00:00:00"
$de.Cells.Item(7, 8).Value = "Include"

$de.Hyperlinks.Add($de.Cells.Item(7, 1), $mdUrl, "", "", $mdName) | Out-Null
$de.Hyperlinks.Add($de.Cells.Item(7, 3), $deXlfUrl, "", "", $deXlfName) | Out-Null
$de.Hyperlinks.Add($de.Cells.Item(8, 1), $cfgUrl, "", "", $cfgName) | Out-Null
